$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 221175
$ws.Range("B4").Value = "Psychiatry"

# "10/09/2025" looks like a date to Excel's auto-detection, so force the
# cell to Text first, enter the value, then restore the default style so
# no stray number-format style is left behind on the cell.
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "10/09/2025"
$ws.Range("C4").Style = "Normal"

$ws.Range("D4").Value = "10:02:55"
$ws.Range("E4").Value = "Scan"
$ws.Range("F4").Value = "160715@med.asu.edu.eg"

# Select the newly-added row, matching the workbook's last on-screen state.
$ws.Rows(4).Select() | Out-Null
